$d = $word.ActiveDocument

# 1. "A list which contains..." -> "It is a list which contains..."
$d.Content.Find.Execute(
    "A list which contains all the time slots",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It is a list which contains all the time slots", 2
) | Out-Null

# 2. Remove the empty paragraph right before the "Fitness:" heading.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "" -and $i -lt $d.Paragraphs.Count) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "Fitness:") {
            $p.Range.Delete() | Out-Null
            break
        }
    }
}

# 3. "In each generation, only two candidates will be remained" (Fitness paragraph)
#    -> "In each generation, only 10% percent of candidates will be remained"
$d.Content.Find.Execute(
    "only two candidates will be remained",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "only 10% percent of candidates will be remained", 2
) | Out-Null

# 4. Crossover paragraph: "only two candidates with highest fitness scores"
#    -> "only 10% percent of candidates with highest fitness scores"
#    The "_GoBack" bookmark (originally collapsed right before "We already know...")
#    ends up re-anchored mid-sentence, right after "...only 10% percent of".
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete() | Out-Null

# Locate the "Crossover:" heading, then its following sentence paragraph,
# so the replace/bookmark placement cannot collide with the similarly
# worded sentence in the Fitness section.
$crossoverHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "Crossover:") {
        $crossoverHeadingIndex = $i
        break
    }
}
$sentencePara = $d.Paragraphs.Item($crossoverHeadingIndex + 1)
$sentenceRange = $sentencePara.Range

$sentenceRange.Find.Execute(
    "only two candidates with highest fitness scores",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "only 10% percent of candidates with highest fitness scores", 2
) | Out-Null

$sentencePara = $d.Paragraphs.Item($crossoverHeadingIndex + 1)
$markRange = $sentencePara.Range
$markRange.Find.Execute(
    "In each generation, only 10% percent of",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0
) | Out-Null
$insertPos = $markRange.End
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 5. Mutation paragraph: describe picking a batch of classes instead of just one.
$d.Content.Find.Execute(
    "randomly pick a course class in its schedule and put it to another random time slot of a random classroom. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "randomly pick a couple of course classes in its schedule and put it to another random time slot of a random classroom. The number of the classes we pick is defined by another parameter called mutation size.",
    2
) | Out-Null
